{"js": "// Split the two M2Doc field-delimiter runs (\"{m\" -> \"{\" + \"m\", and\n// \"'.asImage()}\" -> \"'.asImage()\" + \"}\") into their own <w:r> runs,\n// matching TokenIteratorFieldRewriterSplit's output. We rebuild the\n// whole host paragraph's OOXML (preserving every other run/attribute\n// untouched) and swap it in with insertOoxml(\"Replace\") so the new\n// runs land as genuinely separate <w:r> elements instead of being\n// silently re-coalesced by the host's text-edit APIs.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nparagraphs.items.forEach((p) => p.load(\"text\"));\nawait context.sync();\n\n// Locate the paragraph holding the M2Doc \"{m:'...'.asImage()}\" field.\nconst target = paragraphs.items.filter((p) => p.text.indexOf(\".asImage()\") !== -1)[0];\nif (!target) {\n  throw new Error(\"Could not find the asImage() paragraph\");\n}\n\nconst wNs = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"';\nconst colorRPr = '<w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr>';\n\nconst newParagraphInner =\n  '<w:r ' + wNs + ' w:rsidR=\"00DE6D5A\"><w:t>{</w:t></w:r>' +\n  '<w:r ' + wNs + ' w:rsidR=\"00DE6D5A\"><w:t>m</w:t></w:r>' +\n  '<w:r ' + wNs + '><w:t>:</w:t></w:r>' +\n  '<w:r ' + wNs + ' w:rsidR=\"004B598D\">' + colorRPr + \"<w:t>'</w:t></w:r>\" +\n  '<w:r ' + wNs + ' w:rsidR=\"00071DAB\" w:rsidRPr=\"00071DAB\">' + colorRPr + '<w:t>http</w:t></w:r>' +\n  '<w:r ' + wNs + ' w:rsidR=\"00BC6D60\">' + colorRPr + '<w:t>s</w:t></w:r>' +\n  '<w:bookmarkStart ' + wNs + ' w:id=\"0\" w:name=\"_GoBack\"/>' +\n  '<w:bookmarkEnd ' + wNs + ' w:id=\"0\"/>' +\n  '<w:r ' + wNs + ' w:rsidR=\"00071DAB\" w:rsidRPr=\"00071DAB\">' + colorRPr + '<w:t>://www.m2doc.org/images/logo_M2Doc.png</w:t></w:r>' +\n  '<w:r ' + wNs + ' w:rsidR=\"004B598D\">' + colorRPr + \"<w:t>'.asImage()</w:t></w:r>\" +\n  '<w:r ' + wNs + ' w:rsidR=\"004B598D\">' + colorRPr + '<w:t xml:space=\"preserve\">}</w:t></w:r>';\n\nconst ooxml =\n  '<?xml version=\"1.0\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData><w:document ' + wNs + '><w:body><w:p w:rsidR=\"00C52979\" w:rsidRDefault=\"00C52979\" w:rsidP=\"00F5495F\">' + newParagraphInner + '</w:p></w:body></w:document></pkg:xmlData>' +\n  '</pkg:part></pkg:package>';\n\nconst wholeRange = target.getRange(\"Whole\");\nwholeRange.insertOoxml(ooxml, \"Replace\");\nawait context.sync();\n", "ps1": "# Split the two M2Doc field-delimiter runs (\"{m\" -> \"{\" + \"m\", and\n# \"'.asImage()}\" -> \"'.asImage()\" + \"}\") into their own <w:r> runs,\n# matching TokenIteratorFieldRewriterSplit's output. We rebuild the\n# whole host paragraph's OOXML (preserving every other run/attribute\n# untouched) and swap it in via Range.InsertXML so the new runs land\n# as genuinely separate <w:r> elements instead of being silently\n# re-coalesced by the host's text-edit APIs.\n\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*asImage()*\") {\n        $target = $p\n        break\n    }\n}\nif ($target -eq $null) {\n    throw \"Could not find the asImage() paragraph\"\n}\n\n$wNs = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"'\n$colorRPr = '<w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr>'\n\n$newParagraphInner = (\n    '<w:r ' + $wNs + ' w:rsidR=\"00DE6D5A\"><w:t>{</w:t></w:r>' +\n    '<w:r ' + $wNs + ' w:rsidR=\"00DE6D5A\"><w:t>m</w:t></w:r>' +\n    '<w:r ' + $wNs + '><w:t>:</w:t></w:r>' +\n    '<w:r ' + $wNs + ' w:rsidR=\"004B598D\">' + $colorRPr + \"<w:t>'</w:t></w:r>\" +\n    '<w:r ' + $wNs + ' w:rsidR=\"00071DAB\" w:rsidRPr=\"00071DAB\">' + $colorRPr + '<w:t>http</w:t></w:r>' +\n    '<w:r ' + $wNs + ' w:rsidR=\"00BC6D60\">' + $colorRPr + '<w:t>s</w:t></w:r>' +\n    '<w:bookmarkStart ' + $wNs + ' w:id=\"0\" w:name=\"_GoBack\"/>' +\n    '<w:bookmarkEnd ' + $wNs + ' w:id=\"0\"/>' +\n    '<w:r ' + $wNs + ' w:rsidR=\"00071DAB\" w:rsidRPr=\"00071DAB\">' + $colorRPr + '<w:t>://www.m2doc.org/images/logo_M2Doc.png</w:t></w:r>' +\n    '<w:r ' + $wNs + ' w:rsidR=\"004B598D\">' + $colorRPr + \"<w:t>'.asImage()</w:t></w:r>\" +\n    '<w:r ' + $wNs + ' w:rsidR=\"004B598D\">' + $colorRPr + '<w:t xml:space=\"preserve\">}</w:t></w:r>'\n)\n\n$ooxml = (\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData><w:document ' + $wNs + '><w:body><w:p w:rsidR=\"00C52979\" w:rsidRDefault=\"00C52979\" w:rsidP=\"00F5495F\">' + $newParagraphInner + '</w:p></w:body></w:document></pkg:xmlData>' +\n    '</pkg:part></pkg:package>'\n)\n\n$target.Range.InsertXML($ooxml)\n"}
